$wb = $excel.ActiveWorkbook

# --- Sheet 1: Significant Components ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range('C2').Value = '[''QESL'' ''PPUNIT'' ''QEDLESHI'' ''QHISPC'' ''QNOHLTH'' ''QEXTRCT'' ''QSERV'']'
$ws1.Range('C3').Value = '[''PERCAP'' ''QRICH'' ''MDHSEVAL'']'
$ws1.Range('C4').Value = '[''QRENTER'' ''MEDAGE'' ''QAGEDEP'' ''QSSBEN'']'
$ws1.Range('C6').Value = '[''PPUNIT'' ''QRENTER'' ''QNOAUTO'' ''QPOVTY'']'
$ws1.Range('C7').Value = '[''QAGEDEP'' ''QFEMALE'' ''QFEMLBR'']'

# --- Sheet 2: Loading Factors ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range('A2').Value = 'QESL'
$ws2.Range('A3').Value = 'PPUNIT'
$ws2.Range('A6').Value = 'QNOHLTH'
$ws2.Range('A7').Value = 'QEXTRCT'
$ws2.Range('A9').Value = 'QRICH'
$ws2.Range('A10').Value = 'MDHSEVAL'
$ws2.Range('A15').Value = 'QRENTER'
$ws2.Range('A16').Value = 'QNOAUTO'
$ws2.Range('A17').Value = 'QPOVTY'
$ws2.Range('A18').Value = 'MEDAGE'
$ws2.Range('A19').Value = 'QAGEDEP'
$ws2.Range('A21').Value = 'QFEMALE'
$ws2.Range('A22').Value = 'QFEMLBR'

$ws2.Range('B2').Value = 0.8604417976494784
$ws2.Range('C2').Value = 0.160568581771358
$ws2.Range('D2').Value = -0.06915685180667901
$ws2.Range('E2').Value = 0.09000733896073516
$ws2.Range('F2').Value = 0.1477021107151448
$ws2.Range('G2').Value = -0.1588306519121793
$ws2.Range('B3').Value = 0.6950042396927377
$ws2.Range('C3').Value = -0.03972856865249846
$ws2.Range('D3').Value = -0.1344850090949592
$ws2.Range('E3').Value = 0.2953432226727966
$ws2.Range('F3').Value = -0.5026933097031449
$ws2.Range('G3').Value = 0.05289376696467377
$ws2.Range('B4').Value = 0.8796164256370079
$ws2.Range('C4').Value = 0.201491906632824
$ws2.Range('D4').Value = -0.01904085998538924
$ws2.Range('E4').Value = 0.2532049821778706
$ws2.Range('F4').Value = 0.1445291279046995
$ws2.Range('G4').Value = -0.06380105393169604
$ws2.Range('B5').Value = 0.8323335045057275
$ws2.Range('C5').Value = 0.3030498672351184
$ws2.Range('D5').Value = -0.1374418399146687
$ws2.Range('E5').Value = 0.2467719531841395
$ws2.Range('F5').Value = 0.08043021670754288
$ws2.Range('G5').Value = -0.04898791491532303
$ws2.Range('B6').Value = 0.6263524583266623
$ws2.Range('C6').Value = 0.4231482320438646
$ws2.Range('D6').Value = -0.1289248804973536
$ws2.Range('E6').Value = 0.2413968230640429
$ws2.Range('F6').Value = 0.257759549070146
$ws2.Range('G6').Value = -0.1447924359971534
$ws2.Range('B7').Value = 0.7821909691904856
$ws2.Range('C7').Value = 0.1656086918889093
$ws2.Range('D7').Value = 0.01459395357765911
$ws2.Range('E7').Value = -0.003581474197257717
$ws2.Range('F7').Value = 0.05255974705525795
$ws2.Range('G7').Value = -0.1839728301921103
$ws2.Range('B8').Value = 0.4381629365198772
$ws2.Range('C8').Value = 0.66880820247897
$ws2.Range('D8').Value = -0.2718828411818507
$ws2.Range('E8').Value = 0.324075727470801
$ws2.Range('F8').Value = 0.1684201565279889
$ws2.Range('G8').Value = 0.021689691894979
$ws2.Range('B9').Value = 0.1825438482005984
$ws2.Range('C9').Value = 0.8515610073946787
$ws2.Range('D9').Value = -0.1644290267268402
$ws2.Range('E9').Value = 0.2015465458342662
$ws2.Range('F9').Value = 0.2742661756177662
$ws2.Range('G9').Value = -0.01891201803612595
$ws2.Range('B10').Value = 0.3292221004674073
$ws2.Range('C10').Value = 0.777371739145433
$ws2.Range('D10').Value = -0.04704188475049402
$ws2.Range('E10').Value = 0.2450555885401677
$ws2.Range('F10').Value = -0.02255222619899414
$ws2.Range('G10').Value = -0.00420832700537627
$ws2.Range('B11').Value = 0.1025583099229878
$ws2.Range('C11').Value = 0.2551143919769759
$ws2.Range('D11').Value = -0.01547994434575837
$ws2.Range('E11').Value = 0.4851388598071205
$ws2.Range('F11').Value = 0.07551422482507193
$ws2.Range('G11').Value = 0.04224833066485714
$ws2.Range('B12').Value = 0.1062308591004697
$ws2.Range('C12').Value = 0.2023021799095329
$ws2.Range('D12').Value = -0.08247401214069838
$ws2.Range('E12').Value = 0.5455063355779279
$ws2.Range('F12').Value = 0.3751706954134068
$ws2.Range('G12').Value = -0.02169796957825108
$ws2.Range('B13').Value = 0.4669965697143185
$ws2.Range('C13').Value = 0.3070170503415574
$ws2.Range('D13').Value = -0.1813838185569957
$ws2.Range('E13').Value = 0.3973956104008125
$ws2.Range('F13').Value = 0.3079941395575553
$ws2.Range('G13').Value = -0.01273744352553686
$ws2.Range('B14').Value = 0.3828901314623992
$ws2.Range('C14').Value = 0.1456907906455988
$ws2.Range('D14').Value = -0.02647959372736601
$ws2.Range('E14').Value = 0.7632293195229861
$ws2.Range('F14').Value = 0.03593496885792775
$ws2.Range('G14').Value = 0.1547549998694429
$ws2.Range('B15').Value = 0.04486430913051092
$ws2.Range('C15').Value = 0.2365364242186253
$ws2.Range('D15').Value = -0.4535574168164702
$ws2.Range('E15').Value = 0.00930366945806654
$ws2.Range('F15').Value = 0.7647816228678128
$ws2.Range('G15').Value = -0.06132026733522181
$ws2.Range('B16').Value = 0.1354495595574493
$ws2.Range('C16').Value = 0.04929346000391467
$ws2.Range('D16').Value = -0.07400286519769451
$ws2.Range('E16').Value = 0.2470264182405608
$ws2.Range('F16').Value = 0.6400746721504595
$ws2.Range('G16').Value = 0.03481973803541041
$ws2.Range('B17').Value = 0.4014585206389355
$ws2.Range('C17').Value = 0.1301732394600636
$ws2.Range('D17').Value = -0.3438998606562246
$ws2.Range('E17').Value = 0.2479467978591799
$ws2.Range('F17').Value = 0.4633927662150827
$ws2.Range('G17').Value = 0.04907527643241524
$ws2.Range('B18').Value = -0.2723315417211079
$ws2.Range('C18').Value = -0.2399335003391923
$ws2.Range('D18').Value = 0.789803810423562
$ws2.Range('E18').Value = -0.2212014094941296
$ws2.Range('F18').Value = -0.2219254120386182
$ws2.Range('G18').Value = -0.04090682874677157
$ws2.Range('B19').Value = -0.02669253418912738
$ws2.Range('C19').Value = -0.07422244727089139
$ws2.Range('D19').Value = 0.6301726339865847
$ws2.Range('E19').Value = -0.04529862249890031
$ws2.Range('F19').Value = -0.09030204526530469
$ws2.Range('G19').Value = 0.6657397326359032
$ws2.Range('B20').Value = 0.0161760634816227
$ws2.Range('C20').Value = -0.05131259338132382
$ws2.Range('D20').Value = 0.7963514024883601
$ws2.Range('E20').Value = 0.04548269630222896
$ws2.Range('F20').Value = -0.1206698137313195
$ws2.Range('G20').Value = 0.09809679561554918
$ws2.Range('B21').Value = -0.09613898354356032
$ws2.Range('C21').Value = -0.03199341548043538
$ws2.Range('D21').Value = 0.09767192605087387
$ws2.Range('E21').Value = -0.004977691558137697
$ws2.Range('F21').Value = -0.005878212170707476
$ws2.Range('G21').Value = 0.9521791278791332
$ws2.Range('B22').Value = -0.3572529935701216
$ws2.Range('C22').Value = 0.071125475876885
$ws2.Range('D22').Value = -0.08746849430673782
$ws2.Range('E22').Value = 0.2837560559416707
$ws2.Range('F22').Value = 0.06590328589392641
$ws2.Range('G22').Value = 0.6304728008228307

# --- Sheet 3: All Refactor Variances ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range('B2').Value = 4.648254693095826
$ws3.Range('C2').Value = 2.802762021318809
$ws3.Range('D2').Value = 2.337893633089258
$ws3.Range('E2').Value = 2.20969198795899
$ws3.Range('F2').Value = 2.114856493838247
$ws3.Range('G2').Value = 1.893244857714826
$ws3.Range('H2').Value = 1.032237140243889
$ws3.Range('I2').Value = 4.80776463559952
$ws3.Range('J2').Value = 2.511658534618537
$ws3.Range('K2').Value = 2.204500857808814
$ws3.Range('L2').Value = 2.013880028075578
$ws3.Range('M2').Value = 2.003443446494662
$ws3.Range('N2').Value = 1.882546819392985
$ws3.Range('B3').Value = 0.1721575812257713
$ws3.Range('C3').Value = 0.1038060007895855
$ws3.Range('D3').Value = 0.08658865307737994
$ws3.Range('E3').Value = 0.08184044399848113
$ws3.Range('F3').Value = 0.07832801829030543
$ws3.Range('G3').Value = 0.07012017991536391
$ws3.Range('H3').Value = 0.03823100519421813
$ws3.Range('I3').Value = 0.2289411731237867
$ws3.Range('J3').Value = 0.1196027873627875
$ws3.Range('K3').Value = 0.1049762313242292
$ws3.Range('L3').Value = 0.09589904895597992
$ws3.Range('M3').Value = 0.09540206888069817
$ws3.Range('N3').Value = 0.08964508663776119
$ws3.Range('B4').Value = 0.1721575812257713
$ws3.Range('C4').Value = 0.2759635820153569
$ws3.Range('D4').Value = 0.3625522350927368
$ws3.Range('E4').Value = 0.4443926790912179
$ws3.Range('F4').Value = 0.5227206973815234
$ws3.Range('G4').Value = 0.5928408772968873
$ws3.Range('H4').Value = 0.6310718824911055
$ws3.Range('I4').Value = 0.2289411731237867
$ws3.Range('J4').Value = 0.3485439604865742
$ws3.Range('K4').Value = 0.4535201918108034
$ws3.Range('L4').Value = 0.5494192407667833
$ws3.Range('M4').Value = 0.6448213096474814
$ws3.Range('N4').Value = 0.7344663962852426
$ws3.Range('B5').Value = 0.2728018566540996
$ws3.Range('C5').Value = 0.1644915637499482
$ws3.Range('D5').Value = 0.1372088592120096
$ws3.Range('E5').Value = 0.1296848208090377
$ws3.Range('F5').Value = 0.1241190115793337
$ws3.Range('G5').Value = 0.111112825433721
$ws3.Range('H5').Value = 0.06058106256185002
$ws3.Range('I5').Value = 0.3117108887237279
$ws3.Range('J5').Value = 0.1628431034662853
$ws3.Range('K5').Value = 0.1429285694419434
$ws3.Range('L5').Value = 0.130569689016427
$ws3.Range('M5').Value = 0.1298930344032338
$ws3.Range('N5').Value = 0.1220547149483827

# --- Sheet 4: Final Variances ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range('B2').Value = 4.80776463559952
$ws4.Range('C2').Value = 2.511658534618537
$ws4.Range('D2').Value = 2.204500857808814
$ws4.Range('E2').Value = 2.013880028075578
$ws4.Range('F2').Value = 2.003443446494662
$ws4.Range('G2').Value = 1.882546819392985
$ws4.Range('B3').Value = 0.2289411731237867
$ws4.Range('C3').Value = 0.1196027873627875
$ws4.Range('D3').Value = 0.1049762313242292
$ws4.Range('E3').Value = 0.09589904895597992
$ws4.Range('F3').Value = 0.09540206888069817
$ws4.Range('G3').Value = 0.08964508663776119
$ws4.Range('B4').Value = 0.2289411731237867
$ws4.Range('C4').Value = 0.3485439604865742
$ws4.Range('D4').Value = 0.4535201918108034
$ws4.Range('E4').Value = 0.5494192407667833
$ws4.Range('F4').Value = 0.6448213096474814
$ws4.Range('G4').Value = 0.7344663962852426
$ws4.Range('B5').Value = 0.3117108887237279
$ws4.Range('C5').Value = 0.1628431034662853
$ws4.Range('D5').Value = 0.1429285694419434
$ws4.Range('E5').Value = 0.130569689016427
$ws4.Range('F5').Value = 0.1298930344032338
$ws4.Range('G5').Value = 0.1220547149483827

# --- Sheet 5: Included and Excluded ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range('B2').Value = '[[''QESL'', ''PPUNIT'', ''QEDLESHI'', ''QHISPC'', ''QNOHLTH'', ''QEXTRCT'', ''QSERV'', ''PERCAP'', ''QRICH'', ''MDHSEVAL'', ''QRENTER'', ''MEDAGE'', ''QAGEDEP'', ''QSSBEN'', ''QBLACK'', ''QFAM'', ''QFHH'', ''QNOAUTO'', ''QPOVTY'', ''QFEMALE'', ''QFEMLBR'']]'
